$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell E1: "Toss Winner", matching the style of the other header cells
$ws.Range("E1").Value = "Toss Winner"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New (empty) cell in column E for existing data rows 2 and 3.
# Touching a formatting property (that is a no-op vs. the default format)
# forces the engine to materialize an empty cell record instead of
# leaving the cell completely absent from the sheet.
$ws.Range("E2").Value = ""
$ws.Range("E2").Font.Bold = $false
$ws.Range("E3").Value = ""
$ws.Range("E3").Font.Bold = $false

# New row 4 of data
$ws.Range("A4").Value = "22-03-2025"
$ws.Range("B4").Value = "Kolkata Knight Riders vs Royal Challengers Bengaluru"
$ws.Range("C4").Value = ""
$ws.Range("C4").Font.Bold = $false
$ws.Range("D4").Value = "Royal Challengers Bengaluru"
$ws.Range("E4").Value = "Royal Challengers Bengaluru"
